# Applies the README/docx stats fix-up for the Renaissance fj-kmeans
# (Shenandoah GC, JDK21, heap-8G) table:
#   - rows 1-12 get corrected metric values
#   - rows 44-46 (the tab-separated "raw" rows) get collapsed down to the
#     single headline value that was moved up into rows 1-3.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "1503"
$t.Cell(5,1).Range.Text  = "0.00001"
$t.Cell(6,1).Range.Text  = "0.00059"
$t.Cell(7,1).Range.Text  = "0.00012"
$t.Cell(8,1).Range.Text  = "0.00004"
$t.Cell(9,1).Range.Text  = "0.00017"
$t.Cell(10,1).Range.Text = "0.00018"
$t.Cell(11,1).Range.Text = "0.00022"
$t.Cell(12,1).Range.Text = "0.20185"

$t.Cell(44,1).Range.Text = "99.93"
$t.Cell(45,1).Range.Text = "0.2"
$t.Cell(46,1).Range.Text = "288"
